$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values for case with 380 kV (rows 2-25, columns B,C,D,E,G,I,K,M)
$ws.Cells.Item(2, 2).Value = 10.78335592748847 ; $ws.Cells.Item(2, 3).Value = 12.93183655125447 ; $ws.Cells.Item(2, 4).Value = 6.051680083307279 ; $ws.Cells.Item(2, 5).Value = 10.67501344388617 ; $ws.Cells.Item(2, 7).Value = 3.73922194486029 ; $ws.Cells.Item(2, 9).Value = 48.49108567295852 ; $ws.Cells.Item(2, 11).Value = 13.51171245490564 ; $ws.Cells.Item(2, 13).Value = 15.96349986313935
$ws.Cells.Item(3, 2).Value = 10.8167576846205 ; $ws.Cells.Item(3, 3).Value = 12.55621843885591 ; $ws.Cells.Item(3, 4).Value = 5.945423132733894 ; $ws.Cells.Item(3, 5).Value = 10.41856911483452 ; $ws.Cells.Item(3, 7).Value = 3.743493126148991 ; $ws.Cells.Item(3, 9).Value = 46.989660917696 ; $ws.Cells.Item(3, 11).Value = 13.46730107455785 ; $ws.Cells.Item(3, 13).Value = 15.86615053209072
$ws.Cells.Item(4, 2).Value = 10.84469009016524 ; $ws.Cells.Item(4, 3).Value = 12.32244088705934 ; $ws.Cells.Item(4, 4).Value = 5.881213814876086 ; $ws.Cells.Item(4, 5).Value = 10.2609430477339 ; $ws.Cells.Item(4, 7).Value = 3.746242511741947 ; $ws.Cells.Item(4, 9).Value = 46.04091053947268 ; $ws.Cells.Item(4, 11).Value = 13.44571848631072 ; $ws.Cells.Item(4, 13).Value = 15.81147690813352
$ws.Cells.Item(5, 2).Value = 10.8579193653871 ; $ws.Cells.Item(5, 3).Value = 12.2265235717621 ; $ws.Cells.Item(5, 4).Value = 5.855343966622254 ; $ws.Cells.Item(5, 5).Value = 10.19676196682831 ; $ws.Cells.Item(5, 7).Value = 3.747394959338087 ; $ws.Cells.Item(5, 9).Value = 45.64792008798212 ; $ws.Cells.Item(5, 11).Value = 13.43835681566254 ; $ws.Cells.Item(5, 13).Value = 15.79049563901429
$ws.Cells.Item(6, 2).Value = 10.86022698168385 ; $ws.Cells.Item(6, 3).Value = 12.21056151008728 ; $ws.Cells.Item(6, 4).Value = 5.851067225930084 ; $ws.Cells.Item(6, 5).Value = 10.18611076883242 ; $ws.Cells.Item(6, 7).Value = 3.747588262549566 ; $ws.Cells.Item(6, 9).Value = 45.58229181289096 ; $ws.Cells.Item(6, 11).Value = 13.43722105032703 ; $ws.Cells.Item(6, 13).Value = 15.78709058513542
$ws.Cells.Item(7, 2).Value = 10.84486105591989 ; $ws.Cells.Item(7, 3).Value = 12.32114975742639 ; $ws.Cells.Item(7, 4).Value = 5.880863679044725 ; $ws.Cells.Item(7, 5).Value = 10.26007713770285 ; $ws.Cells.Item(7, 7).Value = 3.746257924081014 ; $ws.Cells.Item(7, 9).Value = 46.03563577326065 ; $ws.Cells.Item(7, 11).Value = 13.44561339733683 ; $ws.Cells.Item(7, 13).Value = 15.81118866914289
$ws.Cells.Item(8, 2).Value = 10.79332056587702 ; $ws.Cells.Item(8, 3).Value = 12.80305793274817 ; $ws.Cells.Item(8, 4).Value = 6.014850138237695 ; $ws.Cells.Item(8, 5).Value = 10.58668252686636 ; $ws.Cells.Item(8, 7).Value = 3.740668408346136 ; $ws.Cells.Item(8, 9).Value = 47.97919455133899 ; $ws.Cells.Item(8, 11).Value = 13.49522207198079 ; $ws.Cells.Item(8, 13).Value = 15.92888603690686
$ws.Cells.Item(9, 2).Value = 10.75195207246949 ; $ws.Cells.Item(9, 3).Value = 13.71740156103533 ; $ws.Cells.Item(9, 4).Value = 6.284194198853861 ; $ws.Cells.Item(9, 5).Value = 11.22180366436555 ; $ws.Cells.Item(9, 7).Value = 3.730707117862968 ; $ws.Cells.Item(9, 9).Value = 51.56190408123409 ; $ws.Cells.Item(9, 11).Value = 13.63735830178261 ; $ws.Cells.Item(9, 13).Value = 16.19931796280952
$ws.Cells.Item(10, 2).Value = 10.75888926384511 ; $ws.Cells.Item(10, 3).Value = 14.3634848447686 ; $ws.Cells.Item(10, 4).Value = 6.483982555446616 ; $ws.Cells.Item(10, 5).Value = 11.68007225313752 ; $ws.Cells.Item(10, 7).Value = 3.72398836288867 ; $ws.Cells.Item(10, 9).Value = 54.03646623571272 ; $ws.Cells.Item(10, 11).Value = 13.76862620093102 ; $ws.Cells.Item(10, 13).Value = 16.42089910652275
$ws.Cells.Item(11, 2).Value = 10.77030012132739 ; $ws.Cells.Item(11, 3).Value = 14.65054445425827 ; $ws.Cells.Item(11, 4).Value = 6.574849110049522 ; $ws.Cells.Item(11, 5).Value = 11.88575348893259 ; $ws.Cells.Item(11, 7).Value = 3.721059987565501 ; $ws.Cells.Item(11, 9).Value = 55.12473842184038 ; $ws.Cells.Item(11, 11).Value = 13.83401916058338 ; $ws.Cells.Item(11, 13).Value = 16.52636392630517
$ws.Cells.Item(12, 2).Value = 10.77581648927783 ; $ws.Cells.Item(12, 3).Value = 14.75816787251038 ; $ws.Cells.Item(12, 4).Value = 6.609220541875986 ; $ws.Cells.Item(12, 5).Value = 11.96316402854382 ; $ws.Cells.Item(12, 7).Value = 3.719969334929413 ; $ws.Cells.Item(12, 9).Value = 55.53122631593153 ; $ws.Cells.Item(12, 11).Value = 13.85958243186244 ; $ws.Cells.Item(12, 13).Value = 16.56694320602776
$ws.Cells.Item(13, 2).Value = 10.77457517614877 ; $ws.Cells.Item(13, 3).Value = 14.73503868657303 ; $ws.Cells.Item(13, 4).Value = 6.601820262434248 ; $ws.Cells.Item(13, 5).Value = 11.94651462048322 ; $ws.Cells.Item(13, 7).Value = 3.720203416775104 ; $ws.Cells.Item(13, 9).Value = 55.44393517582661 ; $ws.Cells.Item(13, 11).Value = 13.85404160447792 ; $ws.Cells.Item(13, 13).Value = 16.55817565306738
$ws.Cells.Item(14, 2).Value = 10.77072996266399 ; $ws.Cells.Item(14, 3).Value = 14.6594208512613 ; $ws.Cells.Item(14, 4).Value = 6.577677821305358 ; $ws.Cells.Item(14, 5).Value = 11.89213206342059 ; $ws.Cells.Item(14, 7).Value = 3.720969893827038 ; $ws.Cells.Item(14, 9).Value = 55.15829425940154 ; $ws.Cells.Item(14, 11).Value = 13.83610632489144 ; $ws.Cells.Item(14, 13).Value = 16.52968971148302
$ws.Cells.Item(15, 2).Value = 10.76853052293324 ; $ws.Cells.Item(15, 3).Value = 14.61295945169068 ; $ws.Cells.Item(15, 4).Value = 6.562883943299062 ; $ws.Cells.Item(15, 5).Value = 11.85875697535289 ; $ws.Cells.Item(15, 7).Value = 3.721441756638005 ; $ws.Cells.Item(15, 9).Value = 54.98259287848868 ; $ws.Cells.Item(15, 11).Value = 13.82522416865155 ; $ws.Cells.Item(15, 13).Value = 16.51232397875605
$ws.Cells.Item(16, 2).Value = 10.75831041874485 ; $ws.Cells.Item(16, 3).Value = 14.34457868976724 ; $ws.Cells.Item(16, 4).Value = 6.478040809551831 ; $ws.Cells.Item(16, 5).Value = 11.66656781621086 ; $ws.Cells.Item(16, 7).Value = 3.724182302414153 ; $ws.Cells.Item(16, 9).Value = 53.96457129441903 ; $ws.Cells.Item(16, 11).Value = 13.76446554291446 ; $ws.Cells.Item(16, 13).Value = 16.41409818743999
$ws.Cells.Item(17, 2).Value = 10.75416177804821 ; $ws.Cells.Item(17, 3).Value = 14.17811135755387 ; $ws.Cells.Item(17, 4).Value = 6.425962186049299 ; $ws.Cells.Item(17, 5).Value = 11.54789647157906 ; $ws.Cells.Item(17, 7).Value = 3.725896221155599 ; $ws.Cells.Item(17, 9).Value = 53.33028723638628 ; $ws.Cells.Item(17, 11).Value = 13.72863558695877 ; $ws.Cells.Item(17, 13).Value = 16.35501486419166
$ws.Cells.Item(18, 2).Value = 10.75255227556576 ; $ws.Cells.Item(18, 3).Value = 14.0817243718552 ; $ws.Cells.Item(18, 4).Value = 6.396008247630592 ; $ws.Cells.Item(18, 5).Value = 11.47938165991867 ; $ws.Cells.Item(18, 7).Value = 3.726894080701893 ; $ws.Cells.Item(18, 9).Value = 52.96195413322424 ; $ws.Cells.Item(18, 11).Value = 13.70856272832593 ; $ws.Cells.Item(18, 13).Value = 16.32147275964805
$ws.Cells.Item(19, 2).Value = 10.75214042603153 ; $ws.Cells.Item(19, 3).Value = 14.04898263747981 ; $ws.Cells.Item(19, 4).Value = 6.385867505973346 ; $ws.Cells.Item(19, 5).Value = 11.45614177334841 ; $ws.Cells.Item(19, 7).Value = 3.7272340144608 ; $ws.Cells.Item(19, 9).Value = 52.83664773673994 ; $ws.Cells.Item(19, 11).Value = 13.70185885261964 ; $ws.Cells.Item(19, 13).Value = 16.3101925918697
$ws.Cells.Item(20, 2).Value = 10.7545229554584 ; $ws.Cells.Item(20, 3).Value = 14.19589898235456 ; $ws.Cells.Item(20, 4).Value = 6.431506271691884 ; $ws.Cells.Item(20, 5).Value = 11.56055652650235 ; $ws.Cells.Item(20, 7).Value = 3.725712524621033 ; $ws.Cells.Item(20, 9).Value = 53.39817289621257 ; $ws.Cells.Item(20, 11).Value = 13.73239442731342 ; $ws.Cells.Item(20, 13).Value = 16.36125894231222
$ws.Cells.Item(21, 2).Value = 10.77182690215041 ; $ws.Cells.Item(21, 3).Value = 14.6816616613509 ; $ws.Cells.Item(21, 4).Value = 6.584770339128246 ; $ws.Cells.Item(21, 5).Value = 11.90811904540581 ; $ws.Cells.Item(21, 7).Value = 3.720744266571607 ; $ws.Cells.Item(21, 9).Value = 55.24234811817035 ; $ws.Cells.Item(21, 11).Value = 13.84135276748906 ; $ws.Cells.Item(21, 13).Value = 16.53803953598315
$ws.Cells.Item(22, 2).Value = 10.79010517262953 ; $ws.Cells.Item(22, 3).Value = 14.99280394435503 ; $ws.Cells.Item(22, 4).Value = 6.684701939955419 ; $ws.Cells.Item(22, 5).Value = 12.13246281749857 ; $ws.Cells.Item(22, 7).Value = 3.71760358363657 ; $ws.Cells.Item(22, 9).Value = 56.41479124260928 ; $ws.Cells.Item(22, 11).Value = 13.91721976165328 ; $ws.Cells.Item(22, 13).Value = 16.65730382455158
$ws.Cells.Item(23, 2).Value = 10.77971002393755 ; $ws.Cells.Item(23, 3).Value = 14.82734987035435 ; $ws.Cells.Item(23, 4).Value = 6.631399320728421 ; $ws.Cells.Item(23, 5).Value = 12.01300684039754 ; $ws.Cells.Item(23, 7).Value = 3.719270142713478 ; $ws.Cells.Item(23, 9).Value = 55.79211157443717 ; $ws.Cells.Item(23, 11).Value = 13.87630787699183 ; $ws.Cells.Item(23, 13).Value = 16.59331903849942
$ws.Cells.Item(24, 2).Value = 10.75435725196941 ; $ws.Cells.Item(24, 3).Value = 14.18785931630884 ; $ws.Cells.Item(24, 4).Value = 6.428999829142161 ; $ws.Cells.Item(24, 5).Value = 11.55483381027399 ; $ws.Cells.Item(24, 7).Value = 3.725795534863882 ; $ws.Cells.Item(24, 9).Value = 53.36749321719781 ; $ws.Cells.Item(24, 11).Value = 13.7306934152 ; $ws.Cells.Item(24, 13).Value = 16.35843466649798
$ws.Cells.Item(25, 2).Value = 10.75664618455339 ; $ws.Cells.Item(25, 3).Value = 13.47408998944043 ; $ws.Cells.Item(25, 4).Value = 6.210844466603636 ; $ws.Cells.Item(25, 5).Value = 11.05110293621663 ; $ws.Cells.Item(25, 7).Value = 3.733295880171582 ; $ws.Cells.Item(25, 9).Value = 50.61917225056929 ; $ws.Cells.Item(25, 11).Value = 13.59414812719113 ; $ws.Cells.Item(25, 13).Value = 16.12203891780445
